$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "asd"
$ws.Range("B3").Value = "sad"
$ws.Range("C3").Value = "sad"
$ws.Range("D3").Value = "sad"
$ws.Range("E3").Value = "participant_19"
$ws.Range("F3").Value = 6615829
